$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '72.416.88'
$ws.Range("E2").Value = '  +4.50%  '

# Row 3
$ws.Range("D3").Value = '4.057.40'
$ws.Range("E3").Value = '  +3.97%  '

# Row 4
$ws.Range("E4").Value = '  -0.02%  '

# Row 5
$ws.Range("D5").Value = "'519.58"
$ws.Range("E5").Value = '  -1.34%  '

# Row 6
$ws.Range("D6").Value = "'147.04"
$ws.Range("E6").Value = '  +1.52%  '

# Row 7
$ws.Range("D7").Value = "'0.729"
$ws.Range("E7").Value = '  +18.78%  '

# Row 8
$ws.Range("D8").Value = '4.055.59'
$ws.Range("E8").Value = '  +4.38%  '

# Row 9
$ws.Range("D9").Value = "'0.999"
$ws.Range("E9").Value = '  +0.09%  '

# Row 10
$ws.Range("D10").Value = "'0.771"
$ws.Range("E10").Value = '  +7.42%  '

# Row 11
$ws.Range("E11").Value = '  +1.41%  '

# Row 12
$ws.Range("E12").Value = '  -2.81%  '

# Row 13
$ws.Range("D13").Value = "'47.85"
$ws.Range("E13").Value = '  +13.76%  '

# Row 14
$ws.Range("D14").Value = "'11.20"
$ws.Range("E14").Value = '  +8.86%  '

# Row 15
$ws.Range("D15").Value = '4.704.02'
$ws.Range("E15").Value = '  +3.79%  '

# Row 16
$ws.Range("D16").Value = '4.078.73'
$ws.Range("E16").Value = '  +4.04%  '

# Row 17
$ws.Range("D17").Value = "'21.37"
$ws.Range("E17").Value = '  +8.20%  '

# Row 18
$ws.Range("D18").Value = "'14.19"
$ws.Range("E18").Value = '  +1.39%  '

# Row 19
$ws.Range("E19").Value = '  -1.11%  '

# Row 20
$ws.Range("E20").Value = '  -1.39%  '

# Row 21
$ws.Range("D21").Value = '72.414.95'
$ws.Range("E21").Value = '  +4.60%  '

# Row 22
$ws.Range("D22").Value = "'444.49"
$ws.Range("E22").Value = '  +4.18%  '

# Row 23
$ws.Range("D23").Value = "'104.76"
$ws.Range("E23").Value = '  +18.55%  '

# Row 24
$ws.Range("D24").Value = "'3.59"
$ws.Range("E24").Value = '  +6.50%  '

# Row 25
$ws.Range("D25").Value = "'14.87"
$ws.Range("E25").Value = '  +5.07%  '

# Row 26
$ws.Range("D26").Value = "'4.01"
$ws.Range("E26").Value = '  -0.94%  '

# Row 27
$ws.Range("D27").Value = "'11.58"
$ws.Range("E27").Value = '  +1.35%  '

# Row 28
$ws.Range("D28").Value = "'11.06"
$ws.Range("E28").Value = '  +4.32%  '

# Row 29
$ws.Range("D29").Value = "'37.80"
$ws.Range("E29").Value = '  +3.71%  '

# Row 30
$ws.Range("E30").Value = '  +2.26%  '

# Row 31
$ws.Range("D31").Value = "'3.25"
$ws.Range("E31").Value = '  +15.18%  '

# Row 32
$ws.Range("E32").Value = '  +4.62%  '

# Row 33
$ws.Range("E33").Value = '  +4.16%  '

# Row 34
$ws.Range("D34").Value = "'681.17"
$ws.Range("E34").Value = '  +0.38%  '

# Row 35
$ws.Range("D35").Value = "'6.85"
$ws.Range("E35").Value = '  +14.49%  '

# Row 36
$ws.Range("D36").Value = "'67.14"
$ws.Range("E36").Value = '  -2.12%  '

# Row 37
$ws.Range("D37").Value = "'43.42"
$ws.Range("E37").Value = '  +8.57%  '

# Row 38
$ws.Range("E38").Value = '  -0.73%  '

# Row 39
$ws.Range("D39").Value = '0.0₃0863'
$ws.Range("E39").Value = '  -2.08%  '

# Row 40
$ws.Range("D40").Value = "'3.54"
$ws.Range("E40").Value = '  +8.77%  '

# Row 41
$ws.Range("D41").Value = "'0.152"
$ws.Range("E41").Value = '  +2.13%  '

# Row 42
$ws.Range("E42").Value = '  +0.04%  '

# Row 43
$ws.Range("D43").Value = "'0.0498"
$ws.Range("E43").Value = '  +3.71%  '

# Row 44
$ws.Range("D44").Value = "'0.998"
$ws.Range("E44").Value = '  -0.25%  '

# Row 45
$ws.Range("B45").Value = 'WEMIXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D45").Value = "'3.26"
$ws.Range("E45").Value = '  +2.67%  '

# Row 46
$ws.Range("B46").Value = 'Stellar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D46").Value = "'0.158"
$ws.Range("E46").Value = '  +12.26%  '

# Row 47
$ws.Range("E47").Value = '  -2.60%  '

# Row 48
$ws.Range("D48").Value = "'3.48"
$ws.Range("E48").Value = '  +3.86%  '

# Row 49
$ws.Range("E49").Value = '  +2.15%  '

# Row 50
$ws.Range("D50").Value = "'9.09"
$ws.Range("E50").Value = '  +6.94%  '

# Row 51
$ws.Range("D51").Value = "'3.34"
$ws.Range("E51").Value = '  +2.45%  '
